$d = $word.ActiveDocument

# --- Simple whole-paragraph text swaps (single-run paragraphs) ---

# Objetivos paragraph: old text -> old "Programa resumido" text
$d.Paragraphs(6).Range.Text = "Concepção de Sistemas de Tratamento de Água em Função da Qualidade da Água Bruta; Projeto de ETAs de Ciclo Completo com Emprego da Decantação ou da Flotação por ar Dissolvido por Clarificação; Projeto de ETAs de Filtração Direta Descendente; Projeto de ETAs de Filtração Direta Ascendente; Projeto de ETAs de Dupla Filtração; Projeto de ETAs por Floto-Filtração; Projeto de ETAs de filtração em Múltiplas Etapas - FiME; Métodos Alternativos de Desinfecção e Adsorção em Carvão Ativado; Tratamento dos Resíduos Gerados nas ETAs e Reuso da Água Recuperada."

# Docente paragraph: old text -> old "Objetivos" text
$d.Paragraphs(8).Range.Text = "Fornecer ao aluno os critérios básicos para elaboração das estações de tratamento de água utilizando diferentes tecnologias de tratamento de água para consumo humano e dos resíduos gerados. Os estudantess irão elaborar projetos utilizando as diversas tecnologias de tratamento."

# "Programa resumido" paragraph: old text -> old "Programa" (long) text
$d.Paragraphs(10).Range.Text = "1 - Concepção de Sistemas de Tratamento de Água em Função da Qualidade da Água Bruta - Tecnologias de Tratamento de Água e dos Resíduos Gerados nos ETAs; 2 - Projeto de ETAs de Ciclo Completo com Emprego da Decantação ou da Flotação por ar Dissolvido para Clarificação; Características de água bruta; características de coagulação e coagulantes; 3 - Projeto de ETAs de Filtração Direta Descendente; características de água bruta; características da coagulação e da floculação; efeito da floculação; filtros com taxa constante e taxa declinante; mecanismo da coagulação e principais coagulantes; 4 - Projeto de ETAs de Filtração Direta Ascendente; características de água bruta; mecanismo da coagulação e principais coagulantes; 5 - Projeto de ETAs de Dupla Filtração; Características de água bruta; mecanismo da coagulação e principais coagulantes; instalação com baterias independentes de filtros ascendentes e descendentes; instalação com filtros ascendentes/descendentes;; 6 - Projeto de ETAs por Floto-Filtração; características de água bruta; características da coagulação e da floculação; características dos filtros; 7 - Projeto de ETAs de Filtração em Múltiplas Etapas - FiME; características de água bruta; instalações de pré-filtração dinâmica; pré-filtração em pedregulho com escoamento ascendente, descendente ou horizontal e filtração lenta em areia; considerações sobre a operação e manutenção; 8 - Métodos Alternativos de Desinfecção e Adsorção em Carvão Ativado; unidades de pré e de pós-desinfecção; características da água e formação sub-produtos; isotermas de adsorção; parâmetros de projeto de adosrção e da câmara de contato; 9 - Tratamento dos resíduos Gerados na ETAs e Reuso da Água Recuperada; tecnologia de tratamento de água e características do sistema de tratamento dos resíduos."

# "Programa" paragraph: old text -> old "Metodo" text
$d.Paragraphs(12).Range.Text = "Aulas expositivas, estudos de projetos sobre as diferentes tecnologias de tratamento; visitas técnicas."

# Bibliografia paragraph (single run w/ internal breaks) -> old "Docente" text
$d.Paragraphs(16).Range.Text = "7455355 - Robson da Silva Rocha"

# --- Rebuild the "Avaliacao" paragraph (multi-run with bold labels) ---

$p = $d.Paragraphs(14)
$r = $p.Range
$pStart = $r.Start
$r.Text = "Método: " + [char]11 + "Avaliação composta por 3 (três) provas, sendo uma substitutiva, e por exercícios sobre as unidades de uma estação de tratamento de água." + [char]11 + "Nota Final = 0,4 x MP + 0,6 x MT" + [char]11 + "MP: média das provas; ME: média de trabalhos " + [char]11 + "* valor mínimo da média das notas das provas (MP) = 5,0" + [char]11 + "* valor mínimo da média das notas dos trabalhos e projetos (MT) = 5,0" + [char]11 + [char]11 + "Critério: " + [char]11 + "Prova única com nota igual ou superior a 5,0." + [char]11 + [char]11 + "Norma de recuperação: " + [char]11 + "DI BERNARDO, L. Métodos e Técnicas de Tratamento de Água. ASSOCIAÇÃO BRASILEIRA DE ENGENHARIA SANITÁRIA E AMBIENTAL & LUIZ DI BERNARDO 2 V., Rio de Janeiro, 1993 (2005)" + [char]11 + "DI BERNARDO, L. Algas e suas Influências na Qualidade da Água e nas Tecnologias de Tratamento ASSOCIAÇÃO BRASILEIRA DE ENGENHARIA SANITÁRIA E AMBIENTAL & LUIZ DI BERNARDO, Rio de Janeiro, 1995." + [char]11 + "PROGRAMA DE PESQUISA EM SANEAMENTO BÁSICO Tratamento de Água de Abastecimento por Filtração em Múltiplas Etapas ASSOCIAÇÃO BRASILEIRA DE ENGENHARIA SANITÁRIA E AMBIENTAL, Rio de Janeiro, 1999 (PROGRAMA DE PESQUISA EM SANEAMENTO BÁSICO. Noções Gerais de Tratamento e Disposição Final de Lodos de Estações de Tratamento de Água ASOCIAÇÃO BRASILEIRA DE ENGENHARIA SANITÁRIA E AMBIENTAL, Rio de Janeiro, 2000."

$d.Range($pStart + 0, $pStart + 8).Font.Bold = $true  # metodo_label
$d.Range($pStart + 352, $pStart + 362).Font.Bold = $true  # criterio_label
$d.Range($pStart + 410, $pStart + 432).Font.Bold = $true  # norma_label
